# The workbook contains one worksheet per year (2000-2100). Each sheet holds
# a small 5x5 "material recycled" table where only cells C2, B4, C4, E4 and
# D5 ever carry a (non-zero) numeric amount. The commit flips the sign of
# every such amount from positive to negative ("+" -> "-") on every sheet,
# leaving zeros, labels and formatting untouched.

$wb = $excel.ActiveWorkbook

$targetCells = @("C2", "B4", "C4", "E4", "D5")

foreach ($ws in $wb.Worksheets) {
    foreach ($addr in $targetCells) {
        $rng = $ws.Range($addr)
        $val = $rng.Value()
        if ($val -ne 0) {
            $rng.Value = (0 - $val)
        }
    }
}
